$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the placeholder numeric values in the "Qty/Item No." column (I10, I14:I18)
# while keeping the existing cell formatting, as part of adding the RFD PO/JO function.
$ws.Range("I10").Value = $null
$ws.Range("I14").Value = $null
$ws.Range("I15").Value = $null
$ws.Range("I16").Value = $null
$ws.Range("I17").Value = $null
$ws.Range("I18").Value = $null

# Move the active selection to the merged I10:K10 range
$ws.Activate()
$ws.Range("I10:K10").Select()
